$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (shifts existing rows 3-9 down to 4-10)
$ws.Rows("3").Insert()

# Populate the newly inserted row 3 with a "wait" step
$ws.Range("B3").Value = "wait"
$ws.Range("C3").Value = "toWait"

# Renumber TESTCASEID column (A2:A10) sequentially TC1..TC9
$ws.Range("A2").Value = "TC1"
$ws.Range("A3").Value = "TC2"
$ws.Range("A4").Value = "TC3"
$ws.Range("A5").Value = "TC4"
$ws.Range("A6").Value = "TC5"
$ws.Range("A7").Value = "TC6"
$ws.Range("A8").Value = "TC7"
$ws.Range("A9").Value = "TC8"
$ws.Range("A10").Value = "TC9"

# Restore the selection to match the post-edit state
$ws.Range("E4").Select()
